$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/volume cells are stored as plain text (e.g. "66.699.37", "  -3.84%  ").
# Some new Price values parse as valid numbers (e.g. "572.61", "1.00"), so a
# leading apostrophe is used to force them to stay text, exactly like typing
# '572.61 into a cell in real Excel, instead of being auto-converted to a
# number (which would also silently drop significant trailing zeros).
$ws.Range("D2").Value = '66.699.37'
$ws.Range("E2").Value = '  -3.84%  '
$ws.Range("D3").Value = '3.313.20'
$ws.Range("E3").Value = '  -0.95%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''572.61'
$ws.Range("E5").Value = '  -3.16%  '
$ws.Range("D6").Value = '''182.66'
$ws.Range("E6").Value = '  -5.27%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '''0.600'
$ws.Range("E8").Value = '  -1.00%  '
$ws.Range("E9").Value = '  -3.29%  '
$ws.Range("E10").Value = '  -1.71%  '
$ws.Range("E11").Value = '  -4.34%  '
$ws.Range("D12").Value = '3.889.51'
$ws.Range("E12").Value = '  -1.00%  '
$ws.Range("E13").Value = '  -0.73%  '
$ws.Range("D14").Value = '''27.18'
$ws.Range("E14").Value = '  -3.54%  '
$ws.Range("D15").Value = '66.740.27'
$ws.Range("E15").Value = '  -3.79%  '
$ws.Range("E16").Value = '  -2.51%  '
$ws.Range("D17").Value = '3.311.67'
$ws.Range("E17").Value = '  -1.75%  '
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").Value = '''13.77'
$ws.Range("E18").Value = '  +0.15%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = '''436.32'
$ws.Range("E19").Value = '  +1.66%  '
$ws.Range("E20").Value = '  -2.31%  '
$ws.Range("D21").Value = '''7.63'
$ws.Range("E21").Value = '  -1.46%  '
$ws.Range("D22").Value = '''73.80'
$ws.Range("E22").Value = '  +0.90%  '
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("E24").Value = '  -0.22%  '
$ws.Range("E25").Value = '  -2.30%  '
$ws.Range("E26").Value = '  +1.27%  '
$ws.Range("D27").Value = '''9.08'
$ws.Range("E27").Value = '  -5.29%  '
$ws.Range("D28").Value = '''1.00'
$ws.Range("E28").Value = '  -2.26%  '
$ws.Range("E29").Value = '  -1.72%  '
$ws.Range("D30").Value = '''22.80'
$ws.Range("E30").Value = '  -0.93%  '
$ws.Range("E31").Value = '  -4.51%  '
$ws.Range("E32").Value = '  +0.12%  '
$ws.Range("D33").Value = '''6.78'
$ws.Range("E33").Value = '  -3.12%  '
$ws.Range("E34").Value = '  -3.76%  '
$ws.Range("E35").Value = '  -0.85%  '
$ws.Range("D36").Value = '''160.25'
$ws.Range("E36").Value = '  -2.77%  '
$ws.Range("E37").Value = '  -3.22%  '
$ws.Range("D38").Value = '''27.31'
$ws.Range("E38").Value = '  +1.11%  '
$ws.Range("D39").Value = '2.812.51'
$ws.Range("E39").Value = '  +2.13%  '
$ws.Range("D40").Value = '''0.791'
$ws.Range("E40").Value = '  -2.38%  '
$ws.Range("E41").Value = '  -2.55%  '
$ws.Range("E42").Value = '  -3.86%  '
$ws.Range("E43").Value = '  -1.42%  '
$ws.Range("D44").Value = '''40.14'
$ws.Range("E44").Value = '  -2.53%  '
$ws.Range("D45").Value = '''24.36'
$ws.Range("E45").Value = '  -3.61%  '
$ws.Range("E46").Value = '  -6.33%  '
$ws.Range("D47").Value = '''319.05'
$ws.Range("E47").Value = '  -7.42%  '
$ws.Range("E48").Value = '  -3.50%  '
$ws.Range("D49").Value = '''0.985'
$ws.Range("E49").Value = '  -1.96%  '
$ws.Range("D50").Value = '''6.18'
$ws.Range("E50").Value = '  -1.56%  '
$ws.Range("D51").Value = '''0.0996'
$ws.Range("E51").Value = '  -1.63%  '
